# Insert a new data row at row 131 ("Feria Lagunitas de Puerto Montt" / "Zapallo"
# sheet). This shifts the previous rows 131-172 down to 132-173 and the new
# row 131 receives fresh data (date serial 44463 = 2021-09-24 / Paine /
# 1a (guarda) / volume 1100).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(131).Insert()

$ws.Cells.Item(131, 1).Value = 4
$ws.Cells.Item(131, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(131, 3).Value = "Los Lagos"
$ws.Cells.Item(131, 4).Value = 44463
$ws.Cells.Item(131, 5).Value = 10
$ws.Cells.Item(131, 6).Value = 100112045
$ws.Cells.Item(131, 7).Value = "Zapallo"
$ws.Cells.Item(131, 8).Value = "Paine"
$ws.Cells.Item(131, 9).Value = "1a (guarda)"
$ws.Cells.Item(131, 10).Value = 1100
$ws.Cells.Item(131, 11).Value = 600
$ws.Cells.Item(131, 12).Value = 600
$ws.Cells.Item(131, 13).Value = 600
$ws.Cells.Item(131, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(131, 15).Value = "Región Metropolitana"
$ws.Cells.Item(131, 16).Value = 600
$ws.Cells.Item(131, 17).Value = 1
$ws.Cells.Item(131, 18).Value = "Hortaliza"
